$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "2025/12/03 01:36"
$ws.Range("B6").Value = "36,672位本"
$ws.Range("C6").Value = "84位 広告・宣伝 (本)"
$ws.Range("D6").Value = "165位商業デザイン"
$ws.Range("E6").Value = "1,839位ビジネス実用本"
$ws.Range("F6").Value = "-"
$ws.Range("G6").Value = "-"
